# The alcohol measurement sheet has a duplicate/unused column M; remove it
# so the old column N ("N") slides left and becomes the new column M.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Columns.Item(13).Delete()

# Selection lands where the deleted column used to be, and the view was
# re-zoomed after the edit.
$ws1.Range("M1").Select()
$excel.ActiveWindow.Zoom = 85

# The other (empty) sheets were re-zoomed too when the workbook was resaved.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$excel.ActiveWindow.Zoom = 85

$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$excel.ActiveWindow.Zoom = 85

# Restore the originally active sheet/tab.
$ws1.Activate()
